# Automatic update of files.
# Increment the "Förändrad" (Changed) date in column C for all data rows
# (rows 2-43) from 2025-04-07 (serial 45754) to 2025-04-08 (serial 45755).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 43 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45754) {
        $cell.Value2 = 45755
    }
}
